# correção das notas do fórum para matc65 em 2021.2
# Zero out all the forum-view grade columns (B..J) for every student row (2..50)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:J50").Value = 0
